$d = $word.ActiveDocument

# 1) "Governor of " -> "Governor of; "
$d.Content.Find.Execute("Governor of ", $true, $false, $false, $false, $false, $true, 1, $false, "Governor of; ", 2)

# 2) subtitle paragraph "Classic & Quantum Mechanics" -> "Classic & Quantum Mechanics / Algorithm"
#    This exact phrase also occurs (a) with a leading space further down, and (b) as the
#    start of "...and all its governing bodies." later on, so scope the Find to just the
#    paragraph whose full text equals "Classic & Quantum Mechanics" and replace only once.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Classic & Quantum Mechanics") {
        $rng = $p.Range
        $rng.Find.Execute("Classic & Quantum Mechanics", $true, $false, $false, $false, $false, $true, 0, $false, "Classic & Quantum Mechanics / Algorithm", 1)
    }
}

# 3) "Sequesterer to ICJ-CIJ on required certification and licensing for;" -> "...for:"
$d.Content.Find.Execute("Sequesterer to ICJ-CIJ on required certification and licensing for;", $true, $false, $false, $false, $false, $true, 1, $false, "Sequesterer to ICJ-CIJ on required certification and licensing for:", 2)

# 4) "Classic & Quantum Mechanics and all its governing bodies." -> "Classic & Quantum Mechanics;  Algorithmic implementation(s)."
$d.Content.Find.Execute("Classic & Quantum Mechanics and all its governing bodies.", $true, $false, $false, $false, $false, $true, 1, $false, "Classic & Quantum Mechanics;  Algorithmic implementation(s).", 2)

# Update paragraph formatting (left indent / first line indent) for the paragraph
# that now contains the "Algorithmic implementation(s)." text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Algorithmic implementation(s).*") {
        $p.Format.LeftIndent = 0
        $p.Format.FirstLineIndent = 0
    }
}
